# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet (holdings detail) right before the
#    "总计" (totals) summary sheet, formatted like the other quarterly
#    sheets (2020-Q4 .. 2021-Q3).
# 2. Prepend a new "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing rows down and renumbering the helper index column.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

function Set-TextRow($ws, $rowRange, $formatDonor, $values) {
    # Write plain-text values (e.g. "16.19") without Excel's smart-entry
    # auto-converting them to numbers, and without leaving a stray
    # "Text" number-format style behind: force text entry via a
    # temporary "@" format, then re-stamp the *formatting only* from an
    # already-unstyled donor cell so the cell ends up with no explicit
    # style (matching the look of the sibling quarter sheets).
    $range = $ws.Range($rowRange)
    $range.NumberFormat = "@"
    $cells = $range.Cells
    for ($i = 0; $i -lt $values.Length; $i++) {
        $cells.Item(1, $i + 1).Value = $values[$i]
    }
    $formatDonor.Copy()
    $range.PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" sheet, positioned after "2021-Q3"
# (i.e. right before "总计"), cloning the look of the existing quarterly
# sheets.
# ---------------------------------------------------------------------

$templateSheet = $wb.Worksheets.Item("2021-Q3")

$newSheet = $wb.Worksheets.Add($null, $templateSheet)
$newSheet.Name = "2022-Q1"

# Bring over the header formatting (bold + border + centered, style used
# by the other quarter sheets) and the numbering-column formatting.
$templateSheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$templateSheet.Range("A2").Copy($newSheet.Range("A2:A6"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# An unstyled cell used purely as a "no formatting" donor for PasteSpecial.
$plainDonor = $templateSheet.Range("B2")

$newSheet.Range("A2").Value = 0
Set-TextRow $newSheet "B2:G2" $plainDonor @("012262", "华宝可持续发展混合A", "16.19", "66.51", "2.13", "0.3448")
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
Set-TextRow $newSheet "B3:G3" $plainDonor @("009989", "华宝研究精选混合", "8.65", "85.40", "3.24", "0.2803")
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
Set-TextRow $newSheet "B4:G4" $plainDonor @("012263", "华宝可持续发展混合C", "6.42", "66.51", "2.13", "0.1367")
$newSheet.Range("H4").Value = 5

$newSheet.Range("A5").Value = 3
Set-TextRow $newSheet "B5:G5" $plainDonor @("004845", "南华瑞盈混合A", "2.90", "82.72", "2.57", "0.0745")
$newSheet.Range("H5").Value = 10

$newSheet.Range("A6").Value = 4
Set-TextRow $newSheet "B6:G6" $plainDonor @("004846", "南华瑞盈混合C", "0.10", "82.72", "2.57", "0.0026")
$newSheet.Range("H6").Value = 10

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" summary sheet, pushing
# the existing rows down by one and renumbering column A (0,1,2,...).
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")

$existing = @()
for ($r = 2; $r -le 5; $r++) {
    $existing += ,@(
        $totalSheet.Cells.Item($r, 2).Value2,
        $totalSheet.Cells.Item($r, 3).Value2,
        $totalSheet.Cells.Item($r, 4).Value2
    )
}

for ($i = 0; $i -lt $existing.Length; $i++) {
    $r = $i + 3
    $row = $existing[$i]
    $totalSheet.Cells.Item($r, 2).Value2 = $row[0]
    $totalSheet.Cells.Item($r, 3).Value2 = $row[1]
    $totalSheet.Cells.Item($r, 4).Value2 = $row[2]
}

# Row 6 is brand-new territory on this sheet; clone the numbering-column
# style from row 2 before overwriting the values.
$totalSheet.Range("A2").Copy($totalSheet.Range("A6"))

for ($r = 2; $r -le 6; $r++) {
    $totalSheet.Cells.Item($r, 1).Value2 = $r - 2
}

$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 5
$totalSheet.Cells.Item(2, 4).Value2 = 0.84

# Adding a sheet makes it the active tab as a side effect; restore the
# workbook's original active sheet so the only visible changes are the
# ones described above.
$wb.Worksheets.Item("2020-Q4").Activate()
